$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data for rows 2-11 (columns A-G)
# A=Colaborador_id, B=Colaborador_nome, C=Departamento, D=Motivo_da_ausência,
# E=Horas_de_ausência, F=Data_da_ausência, G=Salário
$data = @(
    @(79538, "Kaique Viana",               "Vendas",            "Viagem de negócios", 6, 45102, 5650.8),
    @(82832, "Cecília Ferreira",            "Financeiro",        "Outros",              7, 45088, 5030.39),
    @(10735, "Vitor Hugo Jesus",            "Engenharia",        "Problemas pessoais",  2, 45088, 2790.88),
    @(38339, "Gabrielly Fogaça",            "Marketing",         "Viagem de negócios", 8, 45091, 2504.75),
    @(51343, "Alice Silva",                 "P&D",               "Consulta médica",    1, 45102, 4781.13),
    @(12534, "Raquel Moura",                "Marketing",         "Doença",              3, 45078, 12048.55),
    @(1683,  "Sra. Amanda Costa",           "Vendas",            "Outros",              8, 45083, 12080.19),
    @(24829, "Sr. Davi Lucca Cavalcanti",   "Recursos Humanos",  "Doença",              4, 45103, 4722.4),
    @(25815, "Cauê Mendes",                 "Financeiro",        "Problemas pessoais",  8, 45082, 9698.43),
    @(10060, "Clarice da Cunha",            "Recursos Humanos",  "Consulta médica",    6, 45099, 6549.08)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
